$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-32 down to 12-33
$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the new weekly price entry
$ws.Range("A11").Value = 3
$ws.Range("B11").Value = "Femacal de La Calera"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44967
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 100112044
$ws.Range("G11").Value = "Perejil"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 110
$ws.Range("K11").Value = 3000
$ws.Range("L11").Value = 3300
$ws.Range("M11").Value = 3136
$ws.Range("N11").Value = '$/docena de atados (3 kilos)'
$ws.Range("O11").Value = "Provincia de Quillota"
$ws.Range("P11").Value = 1045
$ws.Range("Q11").Value = 3
$ws.Range("R11").Value = "Hortaliza"
